$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ----------------------------------------------------------------------
# Sheet1 ("Runner")
# ----------------------------------------------------------------------

# Existing rows 2-4: Execute column (C) changes from "Yes" to "No"
$ws1.Cells.Item(2,3).Value = "'No"
$ws1.Cells.Item(3,3).Value = "'No"
$ws1.Cells.Item(4,3).Value = "'No"

# New rows 5-7, cloning the formatting from row 4
$ws1.Range("A4:E4").Copy()
$ws1.Range("A5:E7").PasteSpecial(-4122)

$ws1.Cells.Item(5,1).Value = "'bestSellersTest"
$ws1.Cells.Item(5,2).Value = "'To verity header of the page is correct or not"
$ws1.Cells.Item(5,3).Value = "'Yes"
$ws1.Cells.Item(5,4).Value = "'1"
$ws1.Cells.Item(5,5).Value = "'1"

$ws1.Cells.Item(6,1).Value = "'newReleaseTest"
$ws1.Cells.Item(6,2).Value = "'To verity header of the page is correct or not"
$ws1.Cells.Item(6,3).Value = "'Yes"
$ws1.Cells.Item(6,4).Value = "'1"
$ws1.Cells.Item(6,5).Value = "'1"

$ws1.Cells.Item(7,1).Value = "'moversandShakersTest"
$ws1.Cells.Item(7,2).Value = "'To verity header of the page is correct or not"
$ws1.Cells.Item(7,3).Value = "'Yes"
$ws1.Cells.Item(7,4).Value = "'1"
$ws1.Cells.Item(7,5).Value = "'1"

# Column A grows wider to fit the new test names
$ws1.Columns.Item(1).ColumnWidth = 22

# Update selection
[void]$ws1.Range("A5:F7").Select()

# ----------------------------------------------------------------------
# Sheet2 ("CredentialData")
# ----------------------------------------------------------------------
$ws2.Activate()

# Header row: rename "Browser" -> "browser" and add a "menutext" column
$ws2.Cells.Item(1,5).Value = "'browser"
$ws2.Cells.Item(1,5).Copy()
$ws2.Cells.Item(1,6).PasteSpecial(-4122)
$ws2.Cells.Item(1,6).Value = "'menutext"

# Existing rows 2-4: Execute column (D) changes from "Yes" to "No",
# and the browser values move from column E (now re-used) - values stay,
# a new blank "menutext" column F is appended.
$ws2.Cells.Item(2,4).Value = "'No"
$ws2.Cells.Item(2,5).Value = "'chrome"
$ws2.Cells.Item(4,4).Copy()
$ws2.Cells.Item(2,6).PasteSpecial(-4122)
$ws2.Cells.Item(2,6).Value = "'"

$ws2.Cells.Item(3,4).Value = "'No"
$ws2.Cells.Item(3,5).Value = "'firefox"
$ws2.Cells.Item(4,4).Copy()
$ws2.Cells.Item(3,6).PasteSpecial(-4122)
$ws2.Cells.Item(3,6).Value = "'"

$ws2.Cells.Item(4,4).Value = "'No"
$ws2.Cells.Item(4,5).Value = "'edge"
$ws2.Cells.Item(4,4).Copy()
$ws2.Cells.Item(4,6).PasteSpecial(-4122)
$ws2.Cells.Item(4,6).Value = "'"

# New rows 5-7, cloning the formatting from row 4
$ws2.Range("A4:F4").Copy()
$ws2.Range("A5:F7").PasteSpecial(-4122)

$ws2.Cells.Item(5,1).Value = "'bestSellersTest"
$ws2.Cells.Item(5,2).Value = "'"
$ws2.Cells.Item(5,3).Value = "'"
$ws2.Cells.Item(5,4).Value = "'Yes"
$ws2.Cells.Item(5,5).Value = "'"
$ws2.Cells.Item(5,6).Value = "'Best Sellers"

$ws2.Cells.Item(6,1).Value = "'newReleaseTest"
$ws2.Cells.Item(6,2).Value = "'"
$ws2.Cells.Item(6,3).Value = "'"
$ws2.Cells.Item(6,4).Value = "'Yes"
$ws2.Cells.Item(6,5).Value = "'"
$ws2.Cells.Item(6,6).Value = "'New Releases"

$ws2.Cells.Item(7,1).Value = "'moversandShakersTest"
$ws2.Cells.Item(7,2).Value = "'"
$ws2.Cells.Item(7,3).Value = "'"
$ws2.Cells.Item(7,4).Value = "'Yes"
$ws2.Cells.Item(7,5).Value = "'"
$ws2.Cells.Item(7,6).Value = "'Movers and Shakers"

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 22
$ws2.Columns.Item(2).ColumnWidth = 8.5
$ws2.Columns.Item(6).ColumnWidth = 16.59

# Update selection
[void]$ws2.Range("F9").Select()
